$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "candidate" (5th sheet): add party_id column header (H1) first so the
# new shared strings are interned in the same order as the reference edit.
# ---------------------------------------------------------------------------
$wsCandidate = $wb.Worksheets.Item("candidate")

$wsCandidate.Range("H1").Value = "party_id"

# ---------------------------------------------------------------------------
# Sheet "election" (3rd sheet): add a second election row (Europawahl / Wahl 2)
# ---------------------------------------------------------------------------
$wsElection = $wb.Worksheets.Item("election")
$wsElection.Range("A3").Value = 3
$wsElection.Range("B3").Value = 0
$wsElection.Range("C3").Value = "Europawahl"
$wsElection.Range("D3").Value = "Wahl 2"
$wsElection.Range("E3").Value = "2018-05-01T12:00:00.000Z"
$wsElection.Range("F3").Value = "2018-06-29T18:00:00.000Z"
$wsElection.Range("G3").Value = 3
$wsElection.Range("F18").Select()

# ---------------------------------------------------------------------------
# Sheet "candidate" continued: duplicate the candidate rows for the new
# election (id 3) and apply the integer number format to the existing data.
# ---------------------------------------------------------------------------

# give the existing data (D:H, rows 2-25) the integer number format used by
# the new party_id column
$wsCandidate.Range("D2:H25").NumberFormat = "0"

# fill in party_id for the existing AFD/CDU/SPD/FDP candidate rows
$wsCandidate.Range("H22").Value = 50
$wsCandidate.Range("H23").Value = 51
$wsCandidate.Range("H24").Value = 52
$wsCandidate.Range("H25").Value = 53

# duplicate rows 22-25 as rows 26-29, pointing at the new election (id 3)
$candLastNames = @("Gauland", "Kramp-Karrenbauer", "Nahles", "Lindner")
$candFirstNames = @("Alexander", "Annegret", "Andrea", "Christian")
$candParties = @("AFD", "CDU", "SPD", "FDP")
$candD = @(2, 4, 4, 4)
$candG = @(60, 61, 62, 63)
$candH = @(50, 51, 52, 53)

for ($i = 0; $i -lt 4; $i++) {
    $r = 26 + $i
    $wsCandidate.Range("A$r").Value = $candLastNames[$i]
    $wsCandidate.Range("B$r").Value = $candFirstNames[$i]
    $wsCandidate.Range("C$r").Value = $candParties[$i]
    $wsCandidate.Range("D$r").Value = $candD[$i]
    $wsCandidate.Range("E$r").Value = 0
    $wsCandidate.Range("F$r").Value = 3
    $wsCandidate.Range("G$r").Value = $candG[$i]
    $wsCandidate.Range("H$r").Value = $candH[$i]
}

$wsCandidate.Range("D26:H29").NumberFormat = "0"

$wsCandidate.PageSetup.PaperSize = 9
$wsCandidate.PageSetup.Orientation = 1

$wsCandidate.Range("F29").Select()

# ---------------------------------------------------------------------------
# Sheet "party" (6th sheet): duplicate the party rows for the new election
# ---------------------------------------------------------------------------
$wsParty = $wb.Worksheets.Item("party")

$partyNames = @("Alternative für Deutschland", "Christlich Demokratische Union Deutschlands", "Sozialdemokratische Partei Deutschlands", "Freie Demokratische Partei")
$partyC = @(2, 4, 4, 4)
$partyD = @(50, 51, 52, 53)
$partyG = @("AFD", "CDU", "SPD", "FDP")

for ($i = 0; $i -lt 4; $i++) {
    $r = 26 + $i
    $wsParty.Range("B$r").Value = $partyNames[$i]
    $wsParty.Range("C$r").Value = $partyC[$i]
    $wsParty.Range("D$r").Value = $partyD[$i]
    $wsParty.Range("E$r").Value = 3
    $wsParty.Range("F$r").Value = 0
    $wsParty.Range("G$r").Value = $partyG[$i]
}

$wsParty.Range("E29").Select()

# ---------------------------------------------------------------------------
# Sheet "voter" (1st sheet): becomes the active sheet / tab, with a new
# selected cell. Do this last so it ends up as the active tab.
# ---------------------------------------------------------------------------
$wsVoter = $wb.Worksheets.Item("voter")
$wsVoter.Activate()
$wsVoter.Range("O17").Select()
